$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8791316747665405
$ws.Range("B1").Value = 1.740275859832764
$ws.Range("C1").Value = 7.722382545471191
$ws.Range("D1").Value = 2.300868511199951
$ws.Range("E1").Value = 1.500158071517944
